$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data for rows 2-26 (model name + 8 metric columns B..I)
$data = @(
    @("model_9_6_9", 0.7328767345619851, 0.3805099758871867, -3.417986600739741, -0.06971119901748346, 0.2956267297267914, 0.5802342891693115, 0.9144955277442932, 0.7375336885452271),
    @("model_9_6_12", 0.7386749108366446, 0.4688206992886625, -3.713422492051304, -0.04793730382943395, 0.2892098724842072, 0.4975196421146393, 0.9756488800048828, 0.7225213050842285),
    @("model_9_6_8", 0.7400648860367045, 0.3985354511040827, -3.222075901835094, -0.02906873460118553, 0.2876715660095215, 0.563351035118103, 0.8739432096481323, 0.709511935710907),
    @("model_9_6_15", 0.7402085460248129, 0.4955811734983552, -3.856594134331268, -0.04891961712687132, 0.2875125706195831, 0.4724549055099487, 1.005284547805786, 0.7231984734535217),
    @("model_9_6_14", 0.7410105369714249, 0.4970957948519655, -3.833416851156634, -0.04455263387517583, 0.2866249978542328, 0.471036285161972, 1.000486969947815, 0.7201876044273376),
    @("model_9_6_17", 0.7428256317492985, 0.5096458390125644, -3.850825120858295, -0.03798906894783638, 0.2846162617206573, 0.4592815041542053, 1.004090428352356, 0.7156622409820557),
    @("model_9_6_16", 0.7434119287197412, 0.5094229280283432, -3.81288107314812, -0.03278887153165178, 0.2839674055576324, 0.4594902396202087, 0.9962362051010132, 0.7120768427848816),
    @("model_9_6_13", 0.744408237435329, 0.5088959877961279, -3.749958608424032, -0.02427789563681015, 0.2828647792339325, 0.4599838256835938, 0.9832116365432739, 0.7062087059020996),
    @("model_9_6_20", 0.7449117665157929, 0.5341386629614001, -3.90947737287561, -0.02866218765518047, 0.2823074758052826, 0.436340719461441, 1.016231060028076, 0.7092316150665283),
    @("model_9_6_19", 0.7470412308647495, 0.5405011292416837, -3.861162572234061, -0.01725753464813162, 0.2799508273601532, 0.4303814172744751, 1.006230115890503, 0.7013684511184692),
    @("model_9_6_18", 0.7470623542116157, 0.5365870040944307, -3.838490979921124, -0.0168693073790982, 0.2799274325370789, 0.434047520160675, 1.001537322998047, 0.7011008262634277),
    @("model_9_6_21", 0.7485739121947153, 0.5613722332266897, -3.913855054416311, -0.009691270765941118, 0.2782545983791351, 0.4108329117298126, 1.017137169837952, 0.6961517333984375),
    @("model_9_6_10", 0.7508040117480885, 0.486025424555709, -3.326130239702274, 0.01915069638267908, 0.2757865488529205, 0.481405109167099, 0.8954818248748779, 0.676266074180603),
    @("model_9_6_11", 0.7509832823030691, 0.4885617155367659, -3.339314832003774, 0.0191139159367224, 0.2755881249904633, 0.4790295362472534, 0.8982109427452087, 0.6762914061546326),
    @("model_9_6_23", 0.7532162771807592, 0.5766184847328564, -3.819890958157604, 0.01454691396661245, 0.2731168568134308, 0.3965527415275574, 0.9976872205734253, 0.6794402599334717),
    @("model_9_6_24", 0.7532689983366329, 0.5771854746178926, -3.82133647645592, 0.0147517998436365, 0.2730585038661957, 0.396021693944931, 0.997986376285553, 0.6792989373207092),
    @("model_9_6_22", 0.7536555661719005, 0.576820728819007, -3.803862054976987, 0.01695855371040611, 0.2726306915283203, 0.3963633179664612, 0.9943692684173584, 0.6777775287628174),
    @("model_9_6_7", 0.7549227016107316, 0.4295654691935687, -2.763982173971662, 0.05796718551901892, 0.2712283432483673, 0.5342873334884644, 0.7791206836700439, 0.6495032906532288),
    @("model_9_6_6", 0.758311225627606, 0.4380835196485978, -2.673539271157212, 0.07686934589283267, 0.2674782574176788, 0.526309072971344, 0.7603995800018311, 0.6364707946777344),
    @("model_9_6_5", 0.7743732793969738, 0.4757306118273238, -2.204515713690007, 0.1702085912041101, 0.2497023195028305, 0.4910475313663483, 0.6633146405220032, 0.5721161961555481),
    @("model_9_6_4", 0.7986618108137888, 0.5502598722810187, -1.506917313745948, 0.3223696271859453, 0.2228220552206039, 0.4212410748004913, 0.518916130065918, 0.4672057628631592),
    @("model_9_6_3", 0.8503081068627573, 0.7257049924375729, -0.0878365675647772, 0.6490385237456615, 0.1656648069620132, 0.2569135427474976, 0.2251753509044647, 0.2419773489236832),
    @("model_9_6_0", 0.8535807063265186, 0.8346875860686223, 0.719773166631799, 0.841517328554682, 0.1620430052280426, 0.1548369228839874, 0.0580051988363266, 0.1092690378427505),
    @("model_9_6_1", 0.8573265037749034, 0.8274424566445576, 0.6906711651029324, 0.8321955610150702, 0.1578975170850754, 0.1616229563951492, 0.06402912735939026, 0.1156961098313332),
    @("model_9_6_2", 0.8764649190905721, 0.8246128467734545, 0.5020216662531742, 0.8035075973066473, 0.1367169469594955, 0.1642732471227646, 0.1030783951282501, 0.1354756057262421)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $rowNum = $i + 2
    $rowVals = $data[$i]
    $ws.Cells.Item($rowNum, 1).Value = $rowVals[0]
    for ($c = 1; $c -le 8; $c++) {
        $ws.Cells.Item($rowNum, $c + 1).Value = [double]$rowVals[$c]
    }
}
